$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (14) was styled as a plain date (no time); now that
# it's no longer the last row, it reverts to the "date+time" format used by
# all the other interior rows.
$ws.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new trading day: 2021-11-23 (serial 44523) with its total capital.
$ws.Range("A15").Value = 44523
$ws.Range("A15").NumberFormat = "YYYY-MM-DD"
$ws.Range("B15").Value = 60340.15
